# Update the "Pais" COVID tracker sheet with refreshed country data.
# - Two pairs of countries swapped rank order (and picked up new case
#   counts) because the table is kept sorted by "Casos totales" (col B)
#   descending: Emiratos Arabes Unidos overtook Bielorrusia, and
#   Moldavia overtook Argelia.
# - Several other rows just got updated counts without changing rank.
# - The "last updated" timestamp in A1 advanced from 16:52 to 17:22.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Timestamp banner (row 1)
$ws.Range("A1").Value = "Datos actualizados a 24 de Abril de 2020 a las 17:22"

# Row 4: Estados Unidos - refreshed counts, rank unchanged (#1)
$ws.Range("B4").Value = 890198
$ws.Range("C4").Value = 3756
$ws.Range("D4").Value = 89982
$ws.Range("E4").Value = 749813
$ws.Range("G4").Value = 167
$ws.Range("H4").Value = 50403

# Row 18: Suiza - refreshed counts, rank unchanged
$ws.Range("E18").Value = 6499
$ws.Range("G18").Value = 29
$ws.Range("H18").Value = 1578

# Rows 36-37: Emiratos Arabes Unidos overtakes Bielorrusia
$ws.Range("A36").Value = "Emiratos Arabes Unidos"
$ws.Range("B36").Value = 9281
$ws.Range("C36").Value = 525
$ws.Range("D36").Value = 1760
$ws.Range("E36").Value = 7457
$ws.Range("F36").Value = 1
$ws.Range("G36").Value = 8
$ws.Range("H36").Value = 64

$ws.Range("A37").Value = "Bielorrusia"
$ws.Range("B37").Value = 8773
$ws.Range("C37").Value = 751
$ws.Range("D37").Value = 1120
$ws.Range("E37").Value = 7590
$ws.Range("F37").Value = 92
$ws.Range("G37").Value = 3
$ws.Range("H37").Value = 63

# Rows 58-59: Moldavia overtakes Argelia
$ws.Range("A58").Value = "Moldavia"
$ws.Range("B58").Value = 3110
$ws.Range("C58").Value = 184
$ws.Range("D58").Value = 755
$ws.Range("E58").Value = 2271
$ws.Range("F58").Value = 212
$ws.Range("G58").Value = 4
$ws.Range("H58").Value = 84

$ws.Range("A59").Value = "Argelia"
$ws.Range("B59").Value = 3007
$ws.Range("C59").Value = 0
$ws.Range("D59").Value = 1355
$ws.Range("E59").Value = 1245
$ws.Range("F59").Value = 40
$ws.Range("G59").Value = 0
$ws.Range("H59").Value = 407

# Row 91: Republica de Chipre - refreshed counts, rank unchanged
$ws.Range("B91").Value = 804
$ws.Range("C91").Value = 9
$ws.Range("E91").Value = 692
$ws.Range("G91").Value = 1
$ws.Range("H91").Value = 14

# Row 116: Kenia - refreshed counts, rank unchanged
$ws.Range("D116").Value = 94
$ws.Range("E116").Value = 228
